$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking Price strings in column D stay text (matches source formatting).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.219.41'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.682.17'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.06'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5272'
$ws.Range('E6').Value = '  -1.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06364'
$ws.Range('E9').Value = '  -1.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.47'
$ws.Range('E10').Value = '  -2.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07618'
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.689.85'
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.527'
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5758'
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008236'
$ws.Range('E15').Value = '  -2.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.31'
$ws.Range('E16').Value = '  +2.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.237.33'
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.005'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.867'
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('E20').Value = '  -1.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '189.81'
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.231'
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.29'
$ws.Range('E24').Value = '  +2.46%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1260'
$ws.Range('E25').Value = '  -1.40%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.719'
$ws.Range('E26').Value = '  -1.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.82'
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06367'
$ws.Range('E28').Value = '  -1.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.377'
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.315'
$ws.Range('E30').Value = '  -0.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.564'
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('E32').Value = '  -0.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.678'
$ws.Range('E33').Value = '  +0.65%  '
$ws.Range('E34').Value = '  -1.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6125'
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.418'
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.745'
$ws.Range('E37').Value = '  +1.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.176'
$ws.Range('E38').Value = '  -1.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01622'
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.096.37'
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8834'
$ws.Range('E41').Value = '  +1.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.39'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.832.73'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000110'
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.52'
$ws.Range('E46').Value = '  +0.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  -0.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.096'
$ws.Range('E48').Value = '  -0.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05266'
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.004'
$ws.Range('E51').Value = '  -1.26%  '
